$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'44.182.39"
$ws.Range("E2").Value = "  +1.59%  "

$ws.Range("D3").Value = "'2.354.70"
$ws.Range("E3").Value = "  -0.63%  "

$ws.Range("E4").Value = "  +0.21%  "

$ws.Range("E5").Value = "  +4.58%  "

$ws.Range("D6").Value = "'240.65"
$ws.Range("E6").Value = "  +3.06%  "

$ws.Range("D7").Value = "'73.55"
$ws.Range("E7").Value = "  +5.84%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").Value = "'0.558"
$ws.Range("E9").Value = "  +21.75%  "

$ws.Range("D10").Value = "'0.103"
$ws.Range("E10").Value = "  +5.43%  "

$ws.Range("D11").Value = "'31.04"
$ws.Range("E11").Value = "  +17.62%  "

$ws.Range("D12").Value = "'7.47"
$ws.Range("E12").Value = "  +20.04%  "

$ws.Range("E13").Value = "  +2.01%  "

$ws.Range("D14").Value = "'2.705.43"
$ws.Range("E14").Value = "  -0.54%  "

$ws.Range("D15").Value = "'16.79"
$ws.Range("E15").Value = "  +6.97%  "

$ws.Range("D16").Value = "'0.911"
$ws.Range("E16").Value = "  +6.96%  "

$ws.Range("D17").Value = "'2.352.44"
$ws.Range("E17").Value = "  -0.54%  "

$ws.Range("D18").Value = "'44.329.59"
$ws.Range("E18").Value = "  +1.99%  "

$ws.Range("D19").Value = "'0.0000102"
$ws.Range("E19").Value = "  +3.87%  "

$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "'6.63"
$ws.Range("E20").Value = "  +4.63%  "

$ws.Range("B21").Value = "Litecoin"
$ws.Range("C21").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D21").Value = "'77.81"
$ws.Range("E21").Value = "  +5.35%  "

$ws.Range("D22").Value = "'255.59"
$ws.Range("E22").Value = "  +1.85%  "

$ws.Range("E23").Value = "  +0.00%  "

$ws.Range("D24").Value = "'3.79"
$ws.Range("E24").Value = "  -2.95%  "

$ws.Range("E25").Value = "  +3.66%  "

$ws.Range("D26").Value = "'10.64"
$ws.Range("E26").Value = "  +7.00%  "

$ws.Range("E27").Value = "  +3.85%  "

$ws.Range("D28").Value = "'22.63"
$ws.Range("E28").Value = "  -0.41%  "

$ws.Range("D29").Value = "'174.44"
$ws.Range("E29").Value = "  +1.33%  "

$ws.Range("D30").Value = "'1.58"
$ws.Range("E30").Value = "  +1.64%  "

$ws.Range("E31").Value = "  +2.62%  "

$ws.Range("D32").Value = "'0.133"
$ws.Range("E32").Value = "  +4.74%  "

$ws.Range("D33").Value = "'5.39"
$ws.Range("E33").Value = "  +7.85%  "

$ws.Range("E34").Value = "  +8.69%  "

$ws.Range("D35").Value = "'5.37"
$ws.Range("E35").Value = "  +5.84%  "

$ws.Range("D36").Value = "'3.92"
$ws.Range("E36").Value = "  +8.58%  "

$ws.Range("D37").Value = "'2.45"
$ws.Range("E37").Value = "  -0.33%  "

$ws.Range("E38").Value = "  -0.73%  "

$ws.Range("E39").Value = "  +7.53%  "

$ws.Range("D40").Value = "'19.25"
$ws.Range("E40").Value = "  +4.07%  "

$ws.Range("D41").Value = "'9.02"
$ws.Range("E41").Value = "  +0.47%  "

$ws.Range("E42").Value = "  +0.05%  "

$ws.Range("B43").Value = "Cronos"
$ws.Range("C43").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D43").Value = "'0.100"
$ws.Range("E43").Value = "  +5.47%  "

$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").Value = "'1.26"
$ws.Range("E44").Value = "  +3.64%  "

$ws.Range("B45").Value = "Algorand"
$ws.Range("C45").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D45").Value = "'0.187"
$ws.Range("E45").Value = "  +12.27%  "

$ws.Range("D46").Value = "'100.38"
$ws.Range("E46").Value = "  +1.41%  "

$ws.Range("E47").Value = "  -1.12%  "

$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "'2.45"
$ws.Range("E48").Value = "  +9.39%  "

$ws.Range("D49").Value = "'4.49"
$ws.Range("E49").Value = "  +0.16%  "

$ws.Range("D50").Value = "'1.450.55"
$ws.Range("E50").Value = "  +0.06%  "

$ws.Range("B51").Value = "HuobiToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D51").Value = "'2.79"
$ws.Range("E51").Value = "  +2.08%  "
